$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated loading-percent results for the 380 kV case (Case_2_246)
$arr = New-Object 'object[,]' 24,1
$arr[0,0] = 16.01959406949101
$arr[1,0] = 15.47743330537129
$arr[2,0] = 15.13592483686533
$arr[3,0] = 14.99479007538513
$arr[4,0] = 14.97124214392986
$arr[5,0] = 15.13402913170205
$arr[6,0] = 15.8345604890184
$arr[7,0] = 17.13203413288872
$arr[8,0] = 18.02949795789992
$arr[9,0] = 18.42403499875584
$arr[10,0] = 18.57134740879213
$arr[11,0] = 18.53971571731265
$arr[12,0] = 18.43619698878685
$arr[13,0] = 18.37251334877973
$arr[14,0] = 18.00342757140937
$arr[15,0] = 17.77340277991569
$arr[16,0] = 17.63981288901961
$arr[17,0] = 17.59436461695283
$arr[18,0] = 17.79802326911213
$arr[19,0] = 18.46666053672924
$arr[20,0] = 18.89142589096832
$arr[21,0] = 18.66587444435889
$arr[22,0] = 17.78689652810049
$arr[23,0] = 16.79022451877639
$ws.Range("B2:B25").Value = $arr

$arr = New-Object 'object[,]' 24,1
$arr[0,0] = 11.56160852720786
$arr[1,0] = 11.36865213472484
$arr[2,0] = 11.24840691616104
$arr[3,0] = 11.19900722928376
$arr[4,0] = 11.19078173580835
$arr[5,0] = 11.24774224750859
$arr[6,0] = 11.49547159892362
$arr[7,0] = 11.96527893718673
$arr[8,0] = 12.29820290831805
$arr[9,0] = 12.44650785967663
$arr[10,0] = 12.50217650571535
$arr[11,0] = 12.49020971446452
$arr[12,0] = 12.45109781410039
$arr[13,0] = 12.4270755490722
$arr[14,0] = 12.28844417471431
$arr[15,0] = 12.20256458597097
$arr[16,0] = 12.15287521369947
$arr[17,0] = 12.13600202933449
$arr[18,0] = 12.21173730235527
$arr[19,0] = 12.46259956027618
$arr[20,0] = 12.62367269472228
$arr[21,0] = 12.53798081682824
$arr[22,0] = 12.20759129874702
$arr[23,0] = 11.84015998022909
$ws.Range("C2:C25").Value = $arr

$arr = New-Object 'object[,]' 24,1
$arr[0,0] = 12.62012424210234
$arr[1,0] = 12.67513737068431
$arr[2,0] = 12.71070804212539
$arr[3,0] = 12.72565542805511
$arr[4,0] = 12.7281647720885
$arr[5,0] = 12.71090779569355
$arr[6,0] = 12.63872156583369
$arr[7,0] = 12.51132515867677
$arr[8,0] = 12.42627373917748
$arr[9,0] = 12.3894192010721
$arr[10,0] = 12.37572595157918
$arr[11,0] = 12.37866337064734
$arr[12,0] = 12.38828738995186
$arr[13,0] = 12.39421656333312
$arr[14,0] = 12.42871910608677
$arr[15,0] = 12.45035461533854
$arr[16,0] = 12.46297165244823
$arr[17,0] = 12.46727329135672
$arr[18,0] = 12.44803359648165
$arr[19,0] = 12.38545346050756
$arr[20,0] = 12.34608475686551
$arr[21,0] = 12.36695688528367
$arr[22,0] = 12.44908237301484
$arr[23,0] = 12.54428219704282
$ws.Range("E2:E25").Value = $arr

$arr = New-Object 'object[,]' 24,1
$arr[0,0] = 16.86991607391233
$arr[1,0] = 15.89584955866808
$arr[2,0] = 15.26997757108491
$arr[3,0] = 15.00819731993403
$arr[4,0] = 14.96433081551593
$arr[5,0] = 15.26647399323137
$arr[6,0] = 16.53996406344768
$arr[7,0] = 19.00274580682531
$arr[8,0] = 20.67494806633232
$arr[9,0] = 21.3917225636224
$arr[10,0] = 21.65686569030329
$arr[11,0] = 21.60004134736742
$arr[12,0] = 21.4136618050453
$arr[13,0] = 21.29868154950795
$arr[14,0] = 20.62722412089977
$arr[15,0] = 20.20408069597325
$arr[16,0] = 19.95656407809801
$arr[17,0] = 19.87204792380568
$arr[18,0] = 20.24955283636154
$arr[19,0] = 21.46857628470577
$arr[20,0] = 22.22866616901552
$arr[21,0] = 21.82633154458858
$arr[22,0] = 20.22900810905287
$arr[23,0] = 18.34778573295695
$ws.Range("F2:F25").Value = $arr

$arr = New-Object 'object[,]' 24,1
$arr[0,0] = 3.64600382952929
$arr[1,0] = 3.648155361637756
$arr[2,0] = 3.64954555217778
$arr[3,0] = 3.650129508189428
$arr[4,0] = 3.650227528744636
$arr[5,0] = 3.649553356916678
$arr[6,0] = 3.646731361328775
$arr[7,0] = 3.641743480592855
$arr[8,0] = 3.638408176861315
$arr[9,0] = 3.63696159859351
$arr[10,0] = 3.636423921070471
$arr[11,0] = 3.636539270707324
$arr[12,0] = 3.636917161210075
$arr[13,0] = 3.637149945180592
$arr[14,0] = 3.638504131747119
$arr[15,0] = 3.639352944230784
$arr[16,0] = 3.639847813407458
$arr[17,0] = 3.64001651213845
$arr[18,0] = 3.639261898333532
$arr[19,0] = 3.636805891629539
$arr[20,0] = 3.635259655938942
$arr[21,0] = 3.636079538056798
$arr[22,0] = 3.639303038762198
$arr[23,0] = 3.643034749268417
$ws.Range("G2:G25").Value = $arr

$arr = New-Object 'object[,]' 24,1
$arr[0,0] = 9.88586572817043
$arr[1,0] = 9.893599138980672
$arr[2,0] = 9.899691366085174
$arr[3,0] = 9.902512373902034
$arr[4,0] = 9.90300124895205
$arr[5,0] = 9.899728040569942
$arr[6,0] = 9.888253570099993
$arr[7,0] = 9.87639407861052
$arr[8,0] = 9.874137627220891
$arr[9,0] = 9.87450519771606
$arr[10,0] = 9.87484403408439
$arr[11,0] = 9.874762194131234
$arr[12,0] = 9.874529077741318
$arr[13,0] = 9.874412261428891
$arr[14,0] = 9.874141589279041
$arr[15,0] = 9.874332031627647
$arr[16,0] = 9.87457283643846
$arr[17,0] = 9.874676937736542
$arr[18,0] = 9.874298178701006
$arr[19,0] = 9.874592138237851
$arr[20,0] = 9.875947511624576
$arr[21,0] = 9.875117974804247
$arr[22,0] = 9.874313074509299
$arr[23,0] = 9.87846600089526
$ws.Range("L2:L25").Value = $arr

$arr = New-Object 'object[,]' 24,1
$arr[0,0] = 15.00433438426004
$arr[1,0] = 14.88393752972507
$arr[2,0] = 14.81130383484449
$arr[3,0] = 14.78205326360548
$arr[4,0] = 14.77721794452144
$arr[5,0] = 14.81090791043095
$arr[6,0] = 14.96256812926467
$arr[7,0] = 15.26914647680965
$arr[8,0] = 15.4985289045445
$arr[9,0] = 15.60346567392029
$arr[10,0] = 15.64326070915151
$arr[11,0] = 15.63468799663173
$arr[12,0] = 15.60673863228121
$arr[13,0] = 15.5896255610729
$arr[14,0] = 15.49168059183119
$arr[15,0] = 15.43172628081976
$arr[16,0] = 15.39729925391653
$arr[17,0] = 15.38565348670987
$arr[18,0] = 15.43810280516987
$arr[19,0] = 15.61494668540609
$arr[20,0] = 15.73084860004052
$arr[21,0] = 15.66896879738631
$arr[22,0] = 15.43521984715451
$arr[23,0] = 15.18537566566898
$ws.Range("M2:M25").Value = $arr

$arr = New-Object 'object[,]' 24,1
$arr[0,0] = 23.07330189003273
$arr[1,0] = 23.22836449595287
$arr[2,0] = 23.33044211671279
$arr[3,0] = 23.37376212632723
$arr[4,0] = 23.38105928652141
$arr[5,0] = 23.33101937699914
$arr[6,0] = 23.12533911250684
$arr[7,0] = 22.77670904607631
$arr[8,0] = 22.55421468472068
$arr[9,0] = 22.46036757434746
$arr[10,0] = 22.42589528458691
$arr[11,0] = 22.43327201602915
$arr[12,0] = 22.4575101262542
$arr[13,0] = 22.47249562610209
$arr[14,0] = 22.56049658167979
$arr[15,0] = 22.61637350311207
$arr[16,0] = 22.64920543428488
$arr[17,0] = 22.6604405982217
$arr[18,0] = 22.61035354130494
$arr[19,0] = 22.45036183574899
$arr[20,0] = 22.35201200184256
$arr[21,0] = 22.40393247165095
$arr[22,0] = 22.61307296249818
$arr[23,0] = 22.86513618304211
$ws.Range("O2:O25").Value = $arr

